$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.405.03'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '1.721.59'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4906'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.94%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2609'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06187'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '1.727.20'
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07012'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('E12').Value = '  -0.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.558'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5985'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.03%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.19'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = '26.407.01'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.000'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007129'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.89%  '
$ws.Range('D21').Value = '1.942.70'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.472'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('E23').Value = '  -2.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.151'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.28'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.20'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.95%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.396'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.94'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('E29').Value = '  -3.98%  '
$ws.Range('E30').Value = '  -0.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07947'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.666'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.79%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.604'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.54%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9929'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.94%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6232'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9266'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.80%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.391'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.944'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9999'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.11%  '
$ws.Range('E41').Value = '  -1.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.92'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.336'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.77%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3833'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.703'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05362'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.06%  '
$ws.Range('B48').Value = 'Elrond'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.07'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.46%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.680'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.234'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '50.77'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.82%  '
